$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Lesson 41 (row 46): topic, date and video link attached
# ---------------------------------------------------------------------
$ws.Range("C46").Value = "Spring Rest"
$ws.Range("E46").Value = 44260

$ws.Range("F46").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F46"), "https://youtu.be/W7j9-1PfVgg", "", "", "https://youtu.be/W7j9-1PfVgg")
$ws.Range("F46").Value = "https://youtu.be/W7j9-1PfVgg "
$ws.Range("Z1").Copy()
$ws.Range("F46").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 47 held a typo'd lesson number (443); fix it to 42 and attach topic
# ---------------------------------------------------------------------
$ws.Range("B47").Value = 42
$ws.Range("C47").Value = "Spring data jpa, h2 database #1"
$ws.Range("E47").Value = 44263

$ws.Range("F47").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F47"), "https://youtu.be/NSTbUammchE", "", "", "https://youtu.be/NSTbUammchE")
$ws.Range("F47").Value = "https://youtu.be/NSTbUammchE "
$ws.Range("Z1").Copy()
$ws.Range("F47").PasteSpecial(-4122)

$ws.Range("Z1").Clear()

# ---------------------------------------------------------------------
# Row 48: renumber 44 -> 43 and attach its topic
# ---------------------------------------------------------------------
$ws.Range("B48").Value = 43
$ws.Range("C48").Value = "Spring data jpa, h2 database #2"

# ---------------------------------------------------------------------
# Remaining lesson numbers shift up by one now that 41/42/43 are filled
# ---------------------------------------------------------------------
$ws.Range("B49").Value = 44
$ws.Range("B50").Value = 45
$ws.Range("B51").Value = 46
$ws.Range("B52").Value = 47
$ws.Range("B53").Value = 48

# Row 53 becomes the new "current / next lesson" row, carrying the
# highlight style that used to sit on row 54
$ws.Range("B54").Copy()
$ws.Range("B53").PasteSpecial(-4122)

# Row 54 no longer holds a lesson; clear its number and formatting
$ws.Range("B54").Clear()

# ---------------------------------------------------------------------
# Leave the selection where editing left off
# ---------------------------------------------------------------------
$ws.Range("C48").Select()
